# Add a "plotsFile" row to the ProjectConfiguration sheet.
#
# Before:  row 10 = dataFolder, row 11 = dataFile, ... row 14 = outputFolder
# After:   a new row 10 (plotsFile / Plots.xlsx / description) is inserted,
#          pushing the former rows 10-14 down to 11-15.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row above the current row 10 ("dataFolder"), shifting the
# rest of the table down by one row.
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row with the new property.
$ws.Cells.Item(10, 1).Value = "plotsFile"
$ws.Cells.Item(10, 2).Value = "Plots.xlsx"
$ws.Cells.Item(10, 3).Value = 'Name of the excel file with plot definitions. Must be located in the "paramsFolder"'

# Match the selection shown in the edited workbook.
$ws.Range("B10").Select()
